$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Widen column C (ORD/CÓD. area) to fit content - closest achievable value to the
# target stored width of 10.28515625 given this engine's column-width snapping.
$ws.Columns.Item(3).ColumnWidth = 9.5

# Tweak print scale slightly down
$ws.PageSetup.Zoom = 89

# Add a "printed on" stamp to the left of the footer, keep existing right-aligned
# page-number footer
$ws.PageSetup.LeftFooter = "Programação impressa em  &D"
$ws.PageSetup.RightFooter = "Página &P de &N"
